$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: epoch/arrival/collision/local-minima averages after first learning
$ws.Range("A2").Value = 400
$ws.Range("B2").Value = 356
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 35

# Row 5: arrival/collision/minima limit ratios after second learning
$ws.Range("B5").Value = 0.89
$ws.Range("C5").Value = 0.0225
$ws.Range("D5").Value = 0.08749999999999999
